$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# "fixed all list and new tag associations"
#
# List6 / List7 / List8 were erroneously wired up with <w:basedOn w:val="Heading4"/>
# (and the <w:outlineLvl/> + (for List6/List8) <w:contextualSpacing/> that came
# along with being treated like a heading). None of List6/List7/List8 are
# actually referenced anywhere in the document body, so they are rebuilt from
# scratch (Delete + Add) to get a clean paragraph-properties / run-properties
# block with no stray basedOn/outlineLvl left behind, then re-linked to their
# existing linked character styles.
# ---------------------------------------------------------------------------

# --- List6 --------------------------------------------------------------
$old = $d.Styles("List6")
$old.Delete()
$s = $d.Styles.Add("List6", 1)
$s.NameLocal = "List 6"
$s.LinkStyle = $d.Styles("List6Char")
$s.ParagraphFormat.LineSpacingRule = 5
$s.ParagraphFormat.LineSpacing = 12.95
$s.ParagraphFormat.FirstLineIndent = 108
$s.Font.Name = "Times New Roman"
$s.Font.NameBi = "Times New Roman"
$s.Font.Size = 12
$s.Font.SizeBi = 12

# --- List7 ----------------------------------------------------------------
$old = $d.Styles("List7")
$old.Delete()
$s = $d.Styles.Add("List7", 1)
$s.NameLocal = "List 7"
$s.LinkStyle = $d.Styles("List7Char")
$s.NoSpaceBetweenParagraphsOfSameStyle = $true
$s.ParagraphFormat.LineSpacingRule = 5
$s.ParagraphFormat.LineSpacing = 12.95
$s.ParagraphFormat.SpaceAfter = 8
$s.ParagraphFormat.FirstLineIndent = 162
$s.Font.Name = "Times New Roman"
$s.Font.NameBi = "Times New Roman"
$s.Font.Bold = $true
$s.Font.Size = 12
$s.Font.SizeBi = 12

# List7Char picks up the bold that used to live on the List7 paragraph style
$s7c = $d.Styles("List7Char")
$s7c.Font.Bold = $true

# --- List8 ------------------------------------------------------------
$old = $d.Styles("List8")
$old.Delete()
$s = $d.Styles.Add("List8", 1)
$s.NameLocal = "List 8"
$s.LinkStyle = $d.Styles("List8Char")
$s.ParagraphFormat.LineSpacingRule = 5
$s.ParagraphFormat.LineSpacing = 12.95
$s.ParagraphFormat.SpaceAfter = 8
$s.ParagraphFormat.FirstLineIndent = 180
$s.Font.Name = "Times New Roman"
$s.Font.NameBi = "Times New Roman"
$s.Font.Size = 12
$s.Font.SizeBi = 12

# ---------------------------------------------------------------------------
# List1 / List3change / List4change keep their existing basedOn/link, they
# just lose the <w:contextualSpacing/> override. Rebuilt the same way so the
# element is actually removed rather than written out as w:val="0".
# ---------------------------------------------------------------------------

# --- List1 ------------------------------------------------------------
$old = $d.Styles("List1")
$old.Delete()
$s = $d.Styles.Add("List1", 1)
$s.NameLocal = "List 1"
$s.BaseStyle = "List"
$s.LinkStyle = $d.Styles("List1Char")
$s.ParagraphFormat.LineSpacingRule = 0
$s.ParagraphFormat.SpaceAfter = 8
$s.Font.Name = "Times New Roman"
$s.Font.NameBi = "Times New Roman"
$s.Font.Size = 12
$s.Font.SizeBi = 12

# --- List3change --------------------------------------------------------
$old = $d.Styles("List3change")
$old.Delete()
$s = $d.Styles.Add("List3change", 1)
$s.NameLocal = "List 3_change"
$s.BaseStyle = "List1"
$s.LinkStyle = $d.Styles("List3changeChar")
$s.ParagraphFormat.SpaceAfter = 0
$s.ParagraphFormat.LeftIndent = 54
$s.Font.NameFarEast = "Times New Roman"

# --- List4change --------------------------------------------------------
$old = $d.Styles("List4change")
$old.Delete()
$s = $d.Styles.Add("List4change", 1)
$s.NameLocal = "List 4_change"
$s.BaseStyle = "List1"
$s.LinkStyle = $d.Styles("List4changeChar")
$s.ParagraphFormat.SpaceAfter = 0
$s.ParagraphFormat.LeftIndent = 72
$s.Font.NameFarEast = "Times New Roman"
